$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null

$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null
$ws.Range("M3").Value = 2021

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").Value = 2.017314837395458

$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Value = 0.11867182493532386

$ws.Range("L6").Copy() | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").Value = 3.944091449932318

$ws.Range("L7").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null
$ws.Range("M7").Value = 0

$ws.Range("L8").Copy() | Out-Null
$ws.Range("M8").PasteSpecial(-4122) | Out-Null
$ws.Range("M8").Value = "-"

$ws.Range("L9").Copy() | Out-Null
$ws.Range("M9").PasteSpecial(-4122) | Out-Null
$ws.Range("M9").Value = 0

$ws.Range("L10").Copy() | Out-Null
$ws.Range("M10").PasteSpecial(-4122) | Out-Null
$ws.Range("M10").Value = 0.6292103017456653

$ws.Range("L11").Copy() | Out-Null
$ws.Range("M11").PasteSpecial(-4122) | Out-Null
$ws.Range("M11").Value = "-"

$ws.Range("L12").Copy() | Out-Null
$ws.Range("M12").PasteSpecial(-4122) | Out-Null
$ws.Range("M12").Value = 1.2497227177719943

$ws.Range("L13").Copy() | Out-Null
$ws.Range("M13").PasteSpecial(-4122) | Out-Null
$ws.Range("M13").Value = 0.1984453789016842

$ws.Range("L14").Copy() | Out-Null
$ws.Range("M14").PasteSpecial(-4122) | Out-Null
$ws.Range("M14").Value = "-"

$ws.Range("L15").Copy() | Out-Null
$ws.Range("M15").PasteSpecial(-4122) | Out-Null
$ws.Range("M15").Value = 0.39861918314956984

$ws.Range("L16").Copy() | Out-Null
$ws.Range("M16").PasteSpecial(-4122) | Out-Null
$ws.Range("M16").Value = 0

$ws.Range("L17").Copy() | Out-Null
$ws.Range("M17").PasteSpecial(-4122) | Out-Null
$ws.Range("M17").Value = "-"

$ws.Range("L18").Copy() | Out-Null
$ws.Range("M18").PasteSpecial(-4122) | Out-Null
$ws.Range("M18").Value = 0

$ws.Range("L19").Copy() | Out-Null
$ws.Range("M19").PasteSpecial(-4122) | Out-Null
$ws.Range("M19").Value = 0.8552125203112974

$ws.Range("L20").Copy() | Out-Null
$ws.Range("M20").PasteSpecial(-4122) | Out-Null
$ws.Range("M20").Value = "-"

$ws.Range("L21").Copy() | Out-Null
$ws.Range("M21").PasteSpecial(-4122) | Out-Null
$ws.Range("M21").Value = 1.6913581464969858

$ws.Range("L22").Copy() | Out-Null
$ws.Range("M22").PasteSpecial(-4122) | Out-Null
$ws.Range("M22").Value = 1.8347815875998121

$ws.Range("L23").Copy() | Out-Null
$ws.Range("M23").PasteSpecial(-4122) | Out-Null
$ws.Range("M23").Value = "-"

$ws.Range("L24").Copy() | Out-Null
$ws.Range("M24").PasteSpecial(-4122) | Out-Null
$ws.Range("M24").Value = 3.6321107648498847

$ws.Range("L25").Copy() | Out-Null
$ws.Range("M25").PasteSpecial(-4122) | Out-Null
$ws.Range("M25").Value = 6.121156041530003

$ws.Range("L26").Copy() | Out-Null
$ws.Range("M26").PasteSpecial(-4122) | Out-Null
$ws.Range("M26").Value = "-"
$ws.Range("M26").HorizontalAlignment = -4152

$ws.Range("L27").Copy() | Out-Null
$ws.Range("M27").PasteSpecial(-4122) | Out-Null
$ws.Range("M27").Value = 12.437939862560766

$ws.Range("L28").Copy() | Out-Null
$ws.Range("M28").PasteSpecial(-4122) | Out-Null
$ws.Range("M28").Value = 3.6823562661275693

$ws.Range("L29").Copy() | Out-Null
$ws.Range("M29").PasteSpecial(-4122) | Out-Null
$ws.Range("M29").Value = 0.6943323387022582

$ws.Range("L30").Copy() | Out-Null
$ws.Range("M30").PasteSpecial(-4122) | Out-Null
$ws.Range("M30").Value = 7.056499035611798

$ws.Range("L31").Copy() | Out-Null
$ws.Range("M31").PasteSpecial(-4122) | Out-Null
$ws.Range("M31").Value = 2.7447727328177227

$ws.Range("L14").Copy() | Out-Null
$ws.Range("M32").PasteSpecial(-4122) | Out-Null
$ws.Range("M32").Value = "-"

$ws.Range("L33").Copy() | Out-Null
$ws.Range("M33").PasteSpecial(-4122) | Out-Null
$ws.Range("M33").Value = 5.641855041937789

$excel.Application.CutCopyMode = $false
$ws.Range("P6").Select() | Out-Null
